$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row labels
$headers = @("glycan", "binding_score", "monosaccharides", "motifs", "sasa", "flexibility", "has_multi_node_motifs")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Data rows
$data = @(
    @("Gal(b1-4)GlcNAc(b1-2)Man(a1-3)[Gal(b1-4)GlcNAc(b1-2)Man(a1-6)]Man(b1-4)GlcNAc(b1-4)[Fuc(a1-6)]GlcNAc", 9.52787303157406, "['Fuc(a1-6)']", "['Fuc(a1-6)']", 1.656675648076949, 2.802011629266494, $false),
    @("GlcNAc(b1-2)Man(a1-3)[GlcNAc(b1-2)Man(a1-6)]Man(b1-4)GlcNAc(b1-4)[Fuc(a1-6)]GlcNAc", 6.219018257725258, "['Fuc(a1-6)']", "['Fuc(a1-6)']", 2.196111397777991, 2.329919368351718, $false),
    @("Man(a1-2)Man(a1-3)[Man(a1-2)Man(a1-6)]Man(a1-6)[Man(a1-2)Man(a1-2)Man(a1-3)]Man(b1-4)GlcNAc(b1-4)GlcNAc", -0.0996640368447005, "['Man(a1-2)', 'Man(a1-2)', 'Man(a1-2)']", "['Man(a1-2)']", 6.338615601748595, 11.45818660281608, $false),
    @("Man(a1-2)Man(a1-3)[Man(a1-6)]Man(a1-6)[Man(a1-2)Man(a1-3)]Man(b1-4)GlcNAc(b1-4)GlcNAc", -0.088690799271754, "['Man(a1-2)', 'Man(a1-2)']", "['Man(a1-2)']", 4.268688692567416, 5.099384296778321, $false),
    @("Man(a1-2)Man(a1-6)[Man(a1-3)]Man(a1-6)[Man(a1-2)Man(a1-2)Man(a1-3)]Man(b1-4)GlcNAc(b1-4)GlcNAc", -0.1031685426932495, "['Man(a1-2)', 'Man(a1-2)']", "['Man(a1-2)']", 4.967922034550551, 8.378401210032491, $false),
    @("Man(a1-3)[Man(a1-6)]Man(a1-6)[Man(a1-2)Man(a1-3)]Man(b1-4)GlcNAc(b1-4)GlcNAc", 1.050882289654002, "['Man(a1-2)']", "['Man(a1-2)']", 2.414822869487025, 2.445333255766194, $false)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    $excelRow = $r + 2
    $ws.Cells.Item($excelRow, 1).Value = $row[0]
    $ws.Cells.Item($excelRow, 2).Value = $row[1]
    $ws.Cells.Item($excelRow, 3).Value = $row[2]
    $ws.Cells.Item($excelRow, 4).Value = $row[3]
    $ws.Cells.Item($excelRow, 5).Value = $row[4]
    $ws.Cells.Item($excelRow, 6).Value = $row[5]
    $ws.Cells.Item($excelRow, 7).Value = $row[6]
}

# Build the header/label style (bold, thin box border, centered horizontally, top-aligned
# vertically) once on an off-sheet helper cell, then stamp it onto the target cells via
# copy/paste-special so only a single new cell style is materialized.
$helper = $ws.Range("Z100")
$helper.Value = "helper"
$helper.Font.Bold = $true
$helper.Borders.LineStyle = 1
$helper.HorizontalAlignment = -4108
$helper.VerticalAlignment = -4160

$helper.Copy() | Out-Null

$headerRange = $ws.Range("A1:G1")
$headerRange.PasteSpecial(-4122) | Out-Null

$colARange = $ws.Range("A2:A7")
$colARange.PasteSpecial(-4122) | Out-Null

$helper.Delete() | Out-Null

$ws.Range("A1").Select() | Out-Null
